$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Plain value updates (refreshed economic-data pulls / re-run model).
#    These cells keep their existing number formats/styles; only the
#    stored values change.
# ---------------------------------------------------------------------------

$values = @{
    "F7"  = 4.2373

    "F28" = 0.0526514289430049
    "G28" = -0.02093605859677161

    "F29" = 0.1228549628910314
    "G29" = 0.04877483240471108
    "N29" = 46051
    "Q29" = 2.18
    "R29" = 2.22
    "S29" = 2.21
    "T29" = 2.19
    "U29" = 2.18

    "F30" = 0.06490461658514834
    "G30" = -0.0128376635658648
    "N30" = 46051
    "Q30" = 2.35
    "R30" = 2.36
    "S30" = 2.34
    "U30" = 2.32

    "F31" = 0.1254482530098303
    "G31" = 0.05011084527755218

    "C46" = 45962
    "F46" = 292052
    "G46" = 302919
    "H46" = 294225
    "I46" = 284060
    "J46" = 283923

    "C47" = 45962
    "F47" = -0.0358742766218032
    "G47" = 0.02954881468264081
    "H47" = 0.03578469337463908
    "I47" = 0.0004825251916893425
    "J47" = 0.01213465041583639
    "N47" = 46050

    "C48" = 45962
    "F48" = 348877
    "G48" = 332124
    "H48" = 342363
    "I48" = 339690
    "J48" = 358321
    "N48" = 46050
    "Q48" = 3.56
    "R48" = 3.53
    "S48" = 3.56
    "T48" = 3.6
    "U48" = 3.61

    "C49" = 45962
    "F49" = 0.05044200358902096
    "G49" = -0.02990685325224984
    "H49" = 0.007868939327033475
    "I49" = -0.05199527797700942
    "J49" = 0.05791782795597333
    "N49" = 46050
    "Q49" = 3.83
    "R49" = 3.81
    "S49" = 3.82
    "T49" = 3.84
    "U49" = 3.85

    "C50" = 45962
    "F50" = 30075
    "G50" = 29777
    "H50" = 30169
    "I50" = 30416
    "J50" = 28606
    "N50" = 46050
    "Q50" = 4.26
    "R50" = 4.24
    "S50" = 4.22
    "T50" = 4.24

    "C51" = 45962
    "F51" = 0.01000772408234551
    "G51" = -0.01299347011833341
    "H51" = -0.008120725933719042
    "I51" = 0.06327343913864225
    "J51" = 0.02755127698552395
    "N51" = 46048
    "Q51" = 6.1
    "R51" = 6.09
    "S51" = 6.06
    "T51" = 6.16
    "U51" = 6.15

    "N52" = 46050
    "Q52" = 5.88
    "R52" = 5.85
    "S52" = 5.83
    "U52" = 5.85
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# ---------------------------------------------------------------------------
# 2) Move the "most recently refreshed" yellow-highlight formatting.
#    It previously marked the Dur. Orders / ADXDNO date cells (C28:C31);
#    now the Trade Balance / 30y-Mortgage rows (C46:C51, N51) are the
#    freshest pull, so the highlight moves there instead.
#    Use copy/paste-special (formats only) from donor cells that already
#    carry the desired look, so the existing style records are reused
#    instead of creating new duplicate styles.
# ---------------------------------------------------------------------------

# Donor with the plain (non-highlighted) date look.
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C28,C29,C30,C31").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

# Donor with the yellow-highlighted date look.
$ws.Range("N29").Copy() | Out-Null
$ws.Range("C46,C47,C48,C49,C50,C51,N51").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

$excel.CutCopyMode = $false
